$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.490.38"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.98"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.94"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.32"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07892"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9711"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.07"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.833.38"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.888"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.065"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.78"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06639"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001026"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.463.55"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.338"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.036.91"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.50"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.39"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.061"
$ws.Range("E29").Value = "  -4.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.298"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.50"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9432"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09302"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.588"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.255"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.328"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05936"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02187"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.062"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.147"
$ws.Range("E40").Value = "  -4.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5776"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1829"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.994"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.267"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5453"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.94"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.870"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.04"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06605"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.044"
$ws.Range("E51").Value = "  -1.24%  "
